$d = $word.ActiveDocument

# 1) Strike through "Comentar o código-fonte nos quadros"
$found1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Comentar o código-fonte nos quadros*") {
        $found1 = $p
        break
    }
}
$found1.Range.Font.StrikeThrough = 1

# 2) Strike through "Incluir que foi testado no Linux ..." paragraph and move the
#    "_GoBack" bookmark so that it spans this paragraph.
$found2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Incluir que foi testado no Linux*") {
        $found2 = $p
        break
    }
}
$found2.Range.Font.StrikeThrough = 1

# Remove the old "_GoBack" bookmark (currently sitting before "limpar comentários...").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create "_GoBack" so it wraps the "Incluir que foi testado no Linux ..." paragraph.
$d.Bookmarks.Add("_GoBack", $found2.Range)
